# Apply "repull data, push all data, mean calculation" update:
# Column F (dSF) values were refreshed for a subset of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6  = -6
    7  = 7
    8  = -5
    11 = 0
    12 = -2
    15 = -4
    19 = 4
    21 = -12
    22 = 4
    30 = 4
    33 = 13
    37 = -2
    38 = -2
    39 = -4
    40 = -5
    41 = 14
    42 = 2
    43 = 1
    44 = -7
    47 = 0
    49 = -4
    54 = -6
    58 = -2
    59 = -2
    61 = 4
    66 = 1
    71 = 4
    74 = 0
    77 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
